# Updates the "snapshot" sheet's scraped_at column (K) with refreshed timestamps
# and the "returned" sheet's team/player columns (A-D) to match the latest
# publish run, per the 2025-11-24 15:04:54 runtime publish job.
$wb = $excel.ActiveWorkbook
$snap = $wb.Worksheets.Item("snapshot")
$ret = $wb.Worksheets.Item("returned")

# --- snapshot!K2:K38 (scraped_at) ---
$snap.Cells.Item(2, 11).Value = "2025-11-24T07:01:37.931377+00:00"
$snap.Cells.Item(3, 11).Value = "2025-11-24T07:01:40.163933+00:00"
$snap.Cells.Item(4, 11).Value = "2025-11-24T07:01:40.163973+00:00"
$snap.Cells.Item(5, 11).Value = "2025-11-24T07:01:42.870526+00:00"
$snap.Cells.Item(6, 11).Value = "2025-11-24T07:01:45.661702+00:00"
$snap.Cells.Item(7, 11).Value = "2025-11-24T07:01:48.084064+00:00"
$snap.Cells.Item(8, 11).Value = "2025-11-24T07:01:48.084095+00:00"
$snap.Cells.Item(9, 11).Value = "2025-11-24T07:01:48.084114+00:00"
$snap.Cells.Item(10, 11).Value = "2025-11-24T07:01:50.400606+00:00"
$snap.Cells.Item(11, 11).Value = "2025-11-24T07:01:53.177641+00:00"
$snap.Cells.Item(12, 11).Value = "2025-11-24T07:01:55.924296+00:00"
$snap.Cells.Item(13, 11).Value = "2025-11-24T07:01:58.225214+00:00"
$snap.Cells.Item(14, 11).Value = "2025-11-24T07:02:00.541929+00:00"
$snap.Cells.Item(15, 11).Value = "2025-11-24T07:02:05.509866+00:00"
$snap.Cells.Item(16, 11).Value = "2025-11-24T07:02:05.509894+00:00"
$snap.Cells.Item(17, 11).Value = "2025-11-24T07:02:05.509913+00:00"
$snap.Cells.Item(18, 11).Value = "2025-11-24T07:02:07.889722+00:00"
$snap.Cells.Item(19, 11).Value = "2025-11-24T07:02:07.889754+00:00"
$snap.Cells.Item(20, 11).Value = "2025-11-24T07:02:07.889771+00:00"
$snap.Cells.Item(21, 11).Value = "2025-11-24T07:02:10.267375+00:00"
$snap.Cells.Item(22, 11).Value = "2025-11-24T07:02:10.267405+00:00"
$snap.Cells.Item(23, 11).Value = "2025-11-24T07:02:12.975236+00:00"
$snap.Cells.Item(24, 11).Value = "2025-11-24T07:02:12.975263+00:00"
$snap.Cells.Item(25, 11).Value = "2025-11-24T07:02:12.975283+00:00"
$snap.Cells.Item(26, 11).Value = "2025-11-24T07:02:12.975301+00:00"
$snap.Cells.Item(27, 11).Value = "2025-11-24T07:02:15.669380+00:00"
$snap.Cells.Item(28, 11).Value = "2025-11-24T07:02:15.669408+00:00"
$snap.Cells.Item(29, 11).Value = "2025-11-24T07:02:18.399603+00:00"
$snap.Cells.Item(30, 11).Value = "2025-11-24T07:02:18.399637+00:00"
$snap.Cells.Item(31, 11).Value = "2025-11-24T07:02:18.399655+00:00"
$snap.Cells.Item(32, 11).Value = "2025-11-24T07:02:18.399671+00:00"
$snap.Cells.Item(33, 11).Value = "2025-11-24T07:02:20.704368+00:00"
$snap.Cells.Item(34, 11).Value = "2025-11-24T07:02:20.704400+00:00"
$snap.Cells.Item(35, 11).Value = "2025-11-24T07:02:25.876730+00:00"
$snap.Cells.Item(36, 11).Value = "2025-11-24T07:02:25.876759+00:00"
$snap.Cells.Item(37, 11).Value = "2025-11-24T07:02:28.606522+00:00"
$snap.Cells.Item(38, 11).Value = "2025-11-24T07:02:28.606552+00:00"

# --- returned!A3:D37 (team_abbr, team_name, player_name, player_uid) ---
$ret.Cells.Item(3, 1).Value = "СКА"
$ret.Cells.Item(3, 2).Value = "СКА"
$ret.Cells.Item(3, 3).Value = "Короткий Матвей"
$ret.Cells.Item(3, 4).Value = "1369_СКА_короткийматвей"
$ret.Cells.Item(4, 1).Value = "СОЧ"
$ret.Cells.Item(4, 2).Value = "ХК Сочи"
$ret.Cells.Item(4, 3).Value = "Гуськов Матвей"
$ret.Cells.Item(4, 4).Value = "1369_СОЧ_гуськовматвей"
$ret.Cells.Item(5, 1).Value = "СОЧ"
$ret.Cells.Item(5, 2).Value = "ХК Сочи"
$ret.Cells.Item(5, 3).Value = "Самсонов Илья"
$ret.Cells.Item(5, 4).Value = "1369_СОЧ_самсоновилья"
$ret.Cells.Item(6, 1).Value = "СОЧ"
$ret.Cells.Item(6, 2).Value = "ХК Сочи"
$ret.Cells.Item(6, 3).Value = "Сушко Илья"
$ret.Cells.Item(6, 4).Value = "1369_СОЧ_сушкоилья"
$ret.Cells.Item(7, 1).Value = "СОЧ"
$ret.Cells.Item(7, 2).Value = "ХК Сочи"
$ret.Cells.Item(7, 3).Value = "Хёфенмайер Ноэль"
$ret.Cells.Item(7, 4).Value = "1369_СОЧ_хефенмайерноэль"
$ret.Cells.Item(8, 1).Value = "СПР"
$ret.Cells.Item(8, 2).Value = "Спартак"
$ret.Cells.Item(8, 3).Value = "Воробьёв Иван В"
$ret.Cells.Item(8, 4).Value = "1369_СПР_воробьевиванв"
$ret.Cells.Item(9, 1).Value = "СПР"
$ret.Cells.Item(9, 2).Value = "Спартак"
$ret.Cells.Item(9, 3).Value = "Рубцов Герман"
$ret.Cells.Item(9, 4).Value = "1369_СПР_рубцовгерман"
$ret.Cells.Item(10, 1).Value = "СЮЛ"
$ret.Cells.Item(10, 2).Value = "Салават Юлаев"
$ret.Cells.Item(10, 3).Value = "Берлёв Антон"
$ret.Cells.Item(10, 4).Value = "1369_СЮЛ_берлевантон"
$ret.Cells.Item(11, 1).Value = "СЮЛ"
$ret.Cells.Item(11, 2).Value = "Салават Юлаев"
$ret.Cells.Item(11, 3).Value = "Зоркин Никита"
$ret.Cells.Item(11, 4).Value = "1369_СЮЛ_зоркинникита"
$ret.Cells.Item(12, 1).Value = "СЮЛ"
$ret.Cells.Item(12, 2).Value = "Салават Юлаев"
$ret.Cells.Item(12, 3).Value = "Хворов Николай"
$ret.Cells.Item(12, 4).Value = "1369_СЮЛ_хворовниколай"
$ret.Cells.Item(13, 1).Value = "СЮЛ"
$ret.Cells.Item(13, 2).Value = "Салават Юлаев"
$ret.Cells.Item(13, 3).Value = "Ян Денис"
$ret.Cells.Item(13, 4).Value = "1369_СЮЛ_янденис"
$ret.Cells.Item(14, 1).Value = "ТОР"
$ret.Cells.Item(14, 2).Value = "Торпедо"
$ret.Cells.Item(14, 3).Value = "Науменков Михаил"
$ret.Cells.Item(14, 4).Value = "1369_ТОР_науменковмихаил"
$ret.Cells.Item(15, 1).Value = "ТОР"
$ret.Cells.Item(15, 2).Value = "Торпедо"
$ret.Cells.Item(15, 3).Value = "Рожков Никита А"
$ret.Cells.Item(15, 4).Value = "1369_ТОР_рожковникитаа"
$ret.Cells.Item(16, 1).Value = "ЦСК"
$ret.Cells.Item(16, 2).Value = "ЦСКА"
$ret.Cells.Item(16, 3).Value = "Моисеев Данила"
$ret.Cells.Item(16, 4).Value = "1369_ЦСК_моисеевданила"
$ret.Cells.Item(17, 1).Value = "ЦСК"
$ret.Cells.Item(17, 2).Value = "ЦСКА"
$ret.Cells.Item(17, 3).Value = "Уильямс Колби"
$ret.Cells.Item(17, 4).Value = "1369_ЦСК_уильямсколби"
$ret.Cells.Item(18, 1).Value = "СКА"
$ret.Cells.Item(18, 2).Value = "СКА"
$ret.Cells.Item(18, 3).Value = "Зайцев Никита И"
$ret.Cells.Item(18, 4).Value = "1369_СКА_зайцевникитаи"
$ret.Cells.Item(19, 1).Value = "ШДР"
$ret.Cells.Item(19, 2).Value = "Драконы"
$ret.Cells.Item(19, 3).Value = "Бишофф Джейк"
$ret.Cells.Item(19, 4).Value = "1369_ШДР_бишоффджейк"
$ret.Cells.Item(21, 1).Value = "СИБ"
$ret.Cells.Item(21, 2).Value = "Сибирь"
$ret.Cells.Item(21, 3).Value = "Приски Чейз Эванс"
$ret.Cells.Item(21, 4).Value = "1369_СИБ_прискичейзэванс"
$ret.Cells.Item(22, 1).Value = "АВТ"
$ret.Cells.Item(22, 2).Value = "Автомобилист"
$ret.Cells.Item(22, 3).Value = "Трямкин Никита"
$ret.Cells.Item(22, 4).Value = "1369_АВТ_трямкинникита"
$ret.Cells.Item(23, 1).Value = "АДМ"
$ret.Cells.Item(23, 2).Value = "Адмирал"
$ret.Cells.Item(23, 3).Value = "Грман Марио"
$ret.Cells.Item(23, 4).Value = "1369_АДМ_грманмарио"
$ret.Cells.Item(24, 1).Value = "АДМ"
$ret.Cells.Item(24, 2).Value = "Адмирал"
$ret.Cells.Item(24, 3).Value = "Сошников Никита"
$ret.Cells.Item(24, 4).Value = "1369_АДМ_сошниковникита"
$ret.Cells.Item(25, 1).Value = "АКБ"
$ret.Cells.Item(25, 2).Value = "Ак Барс"
$ret.Cells.Item(25, 3).Value = "Яруллин Альберт"
$ret.Cells.Item(25, 4).Value = "1369_АКБ_яруллинальберт"
$ret.Cells.Item(26, 1).Value = "АМР"
$ret.Cells.Item(26, 2).Value = "Амур"
$ret.Cells.Item(26, 3).Value = "Абросимов Роман"
$ret.Cells.Item(26, 4).Value = "1369_АМР_абросимовроман"
$ret.Cells.Item(27, 1).Value = "БАР"
$ret.Cells.Item(27, 2).Value = "Барыс"
$ret.Cells.Item(27, 3).Value = "Бояркин Никита"
$ret.Cells.Item(27, 4).Value = "1369_БАР_бояркинникита"
$ret.Cells.Item(28, 1).Value = "БАР"
$ret.Cells.Item(28, 2).Value = "Барыс"
$ret.Cells.Item(28, 3).Value = "Савицкий Кирилл"
$ret.Cells.Item(28, 4).Value = "1369_БАР_савицкийкирилл"
$ret.Cells.Item(29, 1).Value = "БАР"
$ret.Cells.Item(29, 2).Value = "Барыс"
$ret.Cells.Item(29, 3).Value = "Уотерспун Тайлер"
$ret.Cells.Item(29, 4).Value = "1369_БАР_уотерспунтайлер"
$ret.Cells.Item(30, 1).Value = "ДИН"
$ret.Cells.Item(30, 2).Value = "Динамо М"
$ret.Cells.Item(30, 3).Value = "Готовец Кирилл"
$ret.Cells.Item(30, 4).Value = "1369_ДИН_готовецкирилл"
$ret.Cells.Item(31, 1).Value = "ДМН"
$ret.Cells.Item(31, 2).Value = "Динамо Мн"
$ret.Cells.Item(31, 3).Value = "Уэлле Ксавье"
$ret.Cells.Item(31, 4).Value = "1369_ДМН_уэллексавье"
$ret.Cells.Item(32, 1).Value = "ЛАД"
$ret.Cells.Item(32, 2).Value = "Лада"
$ret.Cells.Item(32, 3).Value = "Ожгихин Алексей"
$ret.Cells.Item(32, 4).Value = "1369_ЛАД_ожгихиналексей"
$ret.Cells.Item(33, 1).Value = "ЛОК"
$ret.Cells.Item(33, 2).Value = "Локомотив"
$ret.Cells.Item(33, 3).Value = "Сергеев Андрей"
$ret.Cells.Item(33, 4).Value = "1369_ЛОК_сергеевандрей"
$ret.Cells.Item(34, 1).Value = "ММГ"
$ret.Cells.Item(34, 2).Value = "Металлург Мг"
$ret.Cells.Item(34, 3).Value = "Козлов Андрей Е"
$ret.Cells.Item(34, 4).Value = "1369_ММГ_козловандрейе"
$ret.Cells.Item(35, 1).Value = "СЕВ"
$ret.Cells.Item(35, 2).Value = "Северсталь"
$ret.Cells.Item(35, 3).Value = "Фомин Макар"
$ret.Cells.Item(35, 4).Value = "1369_СЕВ_фоминмакар"
$ret.Cells.Item(36, 1).Value = "СЕВ"
$ret.Cells.Item(36, 2).Value = "Северсталь"
$ret.Cells.Item(36, 3).Value = "Цицюра Владислав"
$ret.Cells.Item(36, 4).Value = "1369_СЕВ_цицюравладислав"
$ret.Cells.Item(37, 1).Value = "СИБ"
$ret.Cells.Item(37, 2).Value = "Сибирь"
$ret.Cells.Item(37, 3).Value = "Пьянов Валентин"
$ret.Cells.Item(37, 4).Value = "1369_СИБ_пьяноввалентин"
